$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update indicator 3.4.1.1 -> 3.4.1.a (title text of cell B4)
$ws.Range("B4").Value = "3.4.1.a Ожидаемая продолжительность жизни в 15 лет, 45 лет, 65 лет по полу и территории"

# Re-assert the font explicitly so the cell gets its own style record
# (matches how the workbook was actually re-saved after the edit).
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B4").Font.Size = 11

# Move the active selection to B6, matching the saved view state.
$ws.Range("B6").Select()
